$wb = $excel.ActiveWorkbook

# Rename sheets (task order ids refreshed)
$wb.Worksheets.Item("GNG_TO-1650477869264699").Name = "GNG_TO-16509961635872436"
$wb.Worksheets.Item("NB_TO-16504778711436625").Name = "NB_TO-16509961668994176"
$wb.Worksheets.Item("RS_TO-1650477871149665").Name = "RS_TO-16509961668994176"
$wb.Worksheets.Item("TOL_TO-1650477871208665").Name = "TOL_TO-1650996166955416"
$wb.Worksheets.Item("vSAT_TO-1650477871271667").Name = "vSAT_TO-16509961670273814"

# Sheet 1: GNG_TO
$ws1 = $wb.Worksheets.Item("GNG_TO-16509961635872436")
$ws1.Range("B2").Value = "go_stims-16509961635552087.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961635712426.csv"
$ws1.Range("B4").Value = "go_stims-16509961635712426.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961635872436.csv"

# Sheet 2: NB_TO
$ws2 = $wb.Worksheets.Item("NB_TO-16509961668994176")
$ws2.Range("B2").Value = "ZB-match_3-16509961638832119.csv"
$ws2.Range("B3").Value = "OB-16509961642992115.csv"
$ws2.Range("B4").Value = "TB-16509961666833782.csv"
$ws2.Range("B5").Value = "ZB-match_0-1650996163899211.csv"
$ws2.Range("B6").Value = "TB-16509961661553898.csv"
$ws2.Range("B7").Value = "ZB-match_9-16509961638272138.csv"
$ws2.Range("B8").Value = "OB-16509961648594167.csv"
$ws2.Range("B9").Value = "OB-1650996164619382.csv"
$ws2.Range("B10").Value = "TB-16509961668753805.csv"

# Sheet 3: RS_TO
$ws3 = $wb.Worksheets.Item("RS_TO-16509961668994176")
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4: TOL_TO
$ws4 = $wb.Worksheets.Item("TOL_TO-1650996166955416")
$ws4.Range("B2").Value = "MM_stims-16509961669233847.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961668994176.csv"
$ws4.Range("B4").Value = "MM_stims-1650996166939381.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961669233847.csv"
$ws4.Range("B6").Value = "MM_stims-1650996166955416.csv"
$ws4.Range("B7").Value = "ZM_stims-1650996166939381.csv"

# Sheet 5: vSAT_TO
$ws5 = $wb.Worksheets.Item("vSAT_TO-16509961670273814")
$ws5.Range("B2").Value = "vSAT_stims-16509961669874146.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961669714131.csv"
$ws5.Range("B4").Value = "SAT_stims-1650996166955416.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961670114172.csv"
